$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "ДЗ_1"
$ws.Range("D2").Value = "ДЗ_2"
$ws.Range("E2").Value = "ДЗ_3"
$ws.Range("F2").Value = "ДЗ_4"
$ws.Range("G2").Value = "ДЗ_5"

$ws.Range("J8").Select()
